# feat: add 2022-Q4 data
#
# 1. Duplicate the "2022-Q3" sheet (placing the copy right before it),
#    rename the copy to "2022-Q4" and overwrite its fund figures with the
#    new quarter's numbers.
# 2. Insert a new top row in the "总计" (totals) sheet for 2022-Q4 and
#    shift the existing quarters' rows down, renumbering the index column.

$wb = $excel.ActiveWorkbook

# --- 1. New "2022-Q4" detail sheet -----------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.ActiveSheet
$q4.Name = "2022-Q4"

$q4.Range("D2:G2").NumberFormat = "@"
$q4.Cells.Item(2, 4).Value = "0.24"
$q4.Cells.Item(2, 5).Value = "68.42"
$q4.Cells.Item(2, 6).Value = "4.16"
$q4.Cells.Item(2, 7).Value = "0.0100"
$q4.Range("D2:G2").Style = "Normal"
$q4.Cells.Item(2, 8).Value = 8

# --- 2. Update "总计" summary sheet -----------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 1
$total.Cells.Item(2, 4).Value = 0.01

$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(2, 1).PasteSpecial(-4122)

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(6, 1).Value = 4

$wb.Worksheets.Item("总计").Select()
